$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.NotesPage.Shapes.Placeholders.Item(2).TextFrame.TextRange.Text = "Speaker notes here"
